$d = $word.ActiveDocument

# wdFindContinue = 1, wdReplaceAll = 2 (used as literal constants below)

# ---------------------------------------------------------------------------
# Change 1: remove the duplicated " for update UI on screen" trailing text
#   "...4: Door sensor status report for controller update UI on screen
#    for update UI on screen"
#   -> "...4: Door sensor status report for controller update UI on screen"
# ---------------------------------------------------------------------------
$d.Content.Find.Execute(
    "update UI on screen for update UI on screen", $true, $false, $false,
    $false, $false, $true, 1, $false,
    "update UI on screen", 2) | Out-Null

# ---------------------------------------------------------------------------
# Change 2: fix typo "statue" -> "status"
#   "...then update Relay statue and broadcast..."
#   -> "...then update Relay status and broadcast..."
# ---------------------------------------------------------------------------
$d.Content.Find.Execute(
    "update Relay statue and broadcast", $true, $false, $false,
    $false, $false, $true, 1, $false,
    "update Relay status and broadcast", 2) | Out-Null

# ---------------------------------------------------------------------------
# Change 3: insert a new paragraph describing the ESP-NOW "1: Control light
# bulb" behaviour, right after the "When press Switch ON/OFF..." paragraph
# and before the "Door sensor" heading.
# ---------------------------------------------------------------------------
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -eq "When press Switch ON/OFF, it broadcast 2: Light bulb status with broadcast mac address for controller update UI on screen`r") {
        $p.Range.InsertParagraphAfter()
        $newPara = $d.Paragraphs.Item($i + 1)
        $newPara.Range.Text = "When light bulb got 1: Control light bulb from ESP-NOW, it updates Relay status and broadcast 2: Light bulb status"
        break
    }
}

# ---------------------------------------------------------------------------
# Change 4: Door sensor paragraph - add the ESP-NOW "require" clause
#   "if door sensor is power on or door status is change, it broadcast..."
#   -> "if door sensor is power on or door status is change or got
#       5: Require door sensor status from ESP-NOW, it broadcast..."
# ---------------------------------------------------------------------------
$d.Content.Find.Execute(
    "or door status is change, it broadcast", $true, $false, $false,
    $false, $false, $true, 1, $false,
    "or door status is change or got 5: Require door sensor status from ESP-NOW, it broadcast", 2) | Out-Null

# ---------------------------------------------------------------------------
# Change 5: Temperature & Humidity sensor paragraph
#   "...sensor value change, it broadcast 6: Temperature..."
#   -> "...sensor value change or got 7: Require Temperature & Humidity
#       sensor status from ESP-NOW, it broadcast 6: Temperature..."
# ---------------------------------------------------------------------------
$d.Content.Find.Execute(
    "sensor value change, it broadcast 6: Temperature", $true, $false, $false,
    $false, $false, $true, 1, $false,
    "sensor value change or got 7: Require Temperature & Humidity sensor status from ESP-NOW, it broadcast 6: Temperature", 2) | Out-Null

# ---------------------------------------------------------------------------
# Change 6: PM sensor paragraph
#   "if sensor is power on or sensor value change, it broadcast 8: PM..."
#   -> "if sensor is power on or sensor value change or got 9: Require PM
#       sensor status, it broadcast 8: PM..."
# ---------------------------------------------------------------------------
$d.Content.Find.Execute(
    "sensor value change, it broadcast 8: PM", $true, $false, $false,
    $false, $false, $true, 1, $false,
    "sensor value change or got 9: Require PM sensor status, it broadcast 8: PM", 2) | Out-Null
